# Auto-generated Excel COM-interop script to apply Zalera_Profits.xlsx dataset refresh
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) across all 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 4557.6875
$ws.Range("I11").Value = 4557.6875
$ws.Range("K11").Value = 4557.6875
$ws.Range("M11").Value = -4417.6875
$ws.Range("H19").Value = 2401.7693
$ws.Range("I19").Value = 2371.5557
$ws.Range("K19").Value = 2371.5557
$ws.Range("M19").Value = -2196.5557
$ws.Range("H53").Value = 1017.2857
$ws.Range("J53").Value = 3350
$ws.Range("L53").Value = 3350
$ws.Range("N53").Value = -4624
$ws.Range("H138").Value = 3738.6924
$ws.Range("I138").Value = 3999.5
$ws.Range("J138").Value = 3691.2727
$ws.Range("K138").Value = 11998.5
$ws.Range("L138").Value = 11073.8181
$ws.Range("M138").Value = -6858.5
$ws.Range("N138").Value = -21353.8181
$ws.Range("H141").Value = 2706.3333
$ws.Range("I141").Value = 1559.5
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 4678.5
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 501.5
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5151
$ws.Range("I61").Value = 3813.8
$ws.Range("J61").Value = 10499.8
$ws.Range("K61").Value = 3813.8
$ws.Range("L61").Value = 10499.8
$ws.Range("M61").Value = -3601.8
$ws.Range("N61").Value = -10923.8
$ws.Range("H132").Value = 3847.348
$ws.Range("I132").Value = 2934.0588
$ws.Range("K132").Value = 8802.1764
$ws.Range("M132").Value = -6272.1764
$ws.Range("H136").Value = 5151
$ws.Range("I136").Value = 3813.8
$ws.Range("J136").Value = 10499.8
$ws.Range("K136").Value = 11441.4
$ws.Range("L136").Value = 31499.4
$ws.Range("M136").Value = -8891.400000000001
$ws.Range("N136").Value = -36599.39999999999
$ws.Range("H140").Value = 99993.5
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 99993.5
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 99993.5
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -110353.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1358.9286
$ws.Range("I94").Value = 1085.4166
$ws.Range("J94").Value = 3000
$ws.Range("K94").Value = 1085.4166
$ws.Range("L94").Value = 3000
$ws.Range("M94").Value = -634.4166
$ws.Range("N94").Value = -3902
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H103").Value = 250029250
$ws.Range("J103").Value = 250029250
$ws.Range("L103").Value = 250029250
$ws.Range("N103").Value = -250031594
$ws.Range("H105").Value = 100002010
$ws.Range("I105").Value = 111113064
$ws.Range("K105").Value = 111113064
$ws.Range("M105").Value = -111111317
$ws.Range("H134").Value = 5340.125
$ws.Range("I134").Value = 1268.2222
$ws.Range("K134").Value = 3804.6666
$ws.Range("M134").Value = -1269.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H28").Value = 32678.834
$ws.Range("J28").Value = 32678.834
$ws.Range("L28").Value = 32678.834
$ws.Range("N28").Value = -33168.834
$ws.Range("H99").Value = 3960.8
$ws.Range("I99").Value = 3960.8
$ws.Range("K99").Value = 3960.8
$ws.Range("M99").Value = -2462.8
$ws.Range("H122").Value = 1248.5454
$ws.Range("I122").Value = 1341.75
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 4025.25
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -1575.25
$ws.Range("N122").Value = -7900
$ws.Range("H126").Value = 3960.8
$ws.Range("I126").Value = 3960.8
$ws.Range("K126").Value = 11882.4
$ws.Range("M126").Value = -9412.400000000001
$ws.Range("H132").Value = 67801.69
$ws.Range("I132").Value = 3055.1428
$ws.Range("J132").Value = 105570.5
$ws.Range("K132").Value = 9165.428400000001
$ws.Range("L132").Value = 316711.5
$ws.Range("M132").Value = -6635.428400000001
$ws.Range("N132").Value = -321771.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 2745.8076
$ws.Range("I2").Value = 43.666668
$ws.Range("J2").Value = 4176.353
$ws.Range("K2").Value = 262.000008
$ws.Range("L2").Value = 25058.118
$ws.Range("M2").Value = -149.000008
$ws.Range("N2").Value = -25284.118
$ws.Range("H64").Value = 9172.666999999999
$ws.Range("J64").Value = 9008.299999999999
$ws.Range("L64").Value = 27024.9
$ws.Range("N64").Value = -27564.9
$ws.Range("H67").Value = 9172.666999999999
$ws.Range("J67").Value = 9008.299999999999
$ws.Range("L67").Value = 27024.9
$ws.Range("N67").Value = -28896.9
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("H122").Value = 7726.1665
$ws.Range("I122").Value = 15804.8
$ws.Range("J122").Value = 1955.7142
$ws.Range("K122").Value = 142243.2
$ws.Range("L122").Value = 17601.4278
$ws.Range("M122").Value = -139793.2
$ws.Range("N122").Value = -22501.4278
$ws.Range("H131").Value = 16671967
$ws.Range("I131").Value = 83334056
$ws.Range("J131").Value = 6443.8125
$ws.Range("K131").Value = 250002168
$ws.Range("L131").Value = 19331.4375
$ws.Range("M131").Value = -249997128
$ws.Range("N131").Value = -29411.4375
$ws.Range("H138").Value = 3575363.8
$ws.Range("I138").Value = 10000972
$ws.Range("J138").Value = 5581.222
$ws.Range("K138").Value = 30002916
$ws.Range("L138").Value = 16743.666
$ws.Range("M138").Value = -29997776
$ws.Range("N138").Value = -27023.666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3005
$ws.Range("I80").Value = 3005
$ws.Range("K80").Value = 3005
$ws.Range("M80").Value = -2007
$ws.Range("H83").Value = 3005
$ws.Range("I83").Value = 3005
$ws.Range("K83").Value = 15025
$ws.Range("M83").Value = -10033
$ws.Range("H98").Value = 31497.5
$ws.Range("J98").Value = 31497.5
$ws.Range("L98").Value = 31497.5
$ws.Range("N98").Value = -37487.5
$ws.Range("H102").Value = 2693.75
$ws.Range("I102").Value = 2580.8572
$ws.Range("K102").Value = 2580.8572
$ws.Range("M102").Value = -958.8571999999999
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 34633
$ws.Range("I113").Value = 1949.5
$ws.Range("K113").Value = 1949.5
$ws.Range("M113").Value = 220.5
$ws.Range("H122").Value = 8153.6
$ws.Range("I122").Value = 8654.444
$ws.Range("K122").Value = 25963.332
$ws.Range("M122").Value = -23513.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2526.75
$ws.Range("I22").Value = 1490.9678
$ws.Range("K22").Value = 1490.9678
$ws.Range("M22").Value = -1195.9678
$ws.Range("H27").Value = 2526.75
$ws.Range("I27").Value = 1490.9678
$ws.Range("K27").Value = 1490.9678
$ws.Range("M27").Value = -1383.9678
$ws.Range("H100").Value = 19238692
$ws.Range("I100").Value = 31256876
$ws.Range("K100").Value = 31256876
$ws.Range("M100").Value = -31256335
$ws.Range("H132").Value = 5609.8213
$ws.Range("I132").Value = 4382.273
$ws.Range("K132").Value = 13146.819
$ws.Range("M132").Value = -10616.819
$ws.Range("H136").Value = 5354.5884
$ws.Range("I136").Value = 4556.8667
$ws.Range("K136").Value = 13670.6001
$ws.Range("M136").Value = -11120.6001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 17777.777
$ws.Range("I2").Value = 15000
$ws.Range("J2").Value = 27500
$ws.Range("K2").Value = 15000
$ws.Range("L2").Value = 27500
$ws.Range("M2").Value = -14888
$ws.Range("N2").Value = -27724
$ws.Range("H5").Value = 14750
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 14750
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 14750
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -14974
$ws.Range("H136").Value = 7906.5
$ws.Range("I136").Value = 5118.8
$ws.Range("K136").Value = 15356.4
$ws.Range("M136").Value = -12806.4

Write-Host "Applied all Zalera_Profits dataset updates"
